# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-18, replacing the old Strike# derived values
$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
